# ---------------------------------------------------------------------------
# Update BPP_results.xlsx:
#   - run_times: update the three duration cells
#   - scores -> mlp_scores: retitle headers, update MAPE/RMSE values
#   - cat_scores -> mlp_cat_scores: update MAPE/RMSE values (moved to 4th tab)
#   - add new ws_scores sheet (3rd tab) with train/test MAPE & RMSE
#   - add new ws_cat_scores sheet (5th tab, last) with per-instance-size
#     MAPE & RMSE plus a Mean column
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# A cell that already carries the bold / bordered / centred "header" style
# used throughout the workbook (style index 1) - used as a format donor for
# the newly created cells below.
$styleDonor = $wb.Worksheets.Item(1).Range("B1")

# ---------------------------------------------------------------------------
# 1) run_times sheet: just update the three timing strings
# ---------------------------------------------------------------------------
$runTimes = $wb.Worksheets.Item(1)
$runTimes.Range("B2").Value = "7m, 48s"
$runTimes.Range("C2").Value = "11m, 35s"
$runTimes.Range("D2").Value = "0s"

# ---------------------------------------------------------------------------
# 2) scores -> mlp_scores
# ---------------------------------------------------------------------------
$mlpScores = $wb.Worksheets.Item(2)
$mlpScores.Name = "mlp_scores"
$mlpScores.Range("B1").Value = "Train set"
$mlpScores.Range("C1").Value = "Test set"
$mlpScores.Range("B2").Value = 3.19
$mlpScores.Range("C2").Value = 3.03
$mlpScores.Range("B3").Value = 0.01
$mlpScores.Range("C3").Value = 0.01

# ---------------------------------------------------------------------------
# 3) new ws_scores sheet: the old "cat_scores" tab (3rd sheet, sheetId 3) is
#    repurposed in place (cleared + rewritten) so the sheetId sequence stays
#    1,2,3,4,5 and the tab ends up 3rd, exactly where the diff puts it.
# ---------------------------------------------------------------------------
$wsScores = $wb.Worksheets.Item(3)
$wsScores.Cells.Clear()

$wsScores.Range("B1").Value = "Train set"
$wsScores.Range("C1").Value = "Test set"
$wsScores.Range("A2").Value = "MAPE"
$wsScores.Range("B2").Value = 22.21
$wsScores.Range("C2").Value = 22
$wsScores.Range("A3").Value = "RMSE"
$wsScores.Range("B3").Value = 0.09
$wsScores.Range("C3").Value = 0.09

$styleDonor.Copy() | Out-Null
$wsScores.Range("B1:C1").PasteSpecial($xlPasteFormats) | Out-Null
$wsScores.Range("A2").PasteSpecial($xlPasteFormats) | Out-Null
$wsScores.Range("A3").PasteSpecial($xlPasteFormats) | Out-Null

$wsScores.Name = "ws_scores"

# ---------------------------------------------------------------------------
# 4) new mlp_cat_scores sheet (sheetId 4), inserted right after ws_scores -
#    this carries the former "cat_scores" per-instance-size values, updated.
# ---------------------------------------------------------------------------
$mlpCatScores = $wb.Worksheets.Add($null, $wsScores)
$mlpCatScores.Name = "mlp_cat_scores"

$mlpCatScores.Range("B1").Value = 7
$mlpCatScores.Range("C1").Value = 8
$mlpCatScores.Range("D1").Value = 9
$mlpCatScores.Range("E1").Value = 10
$mlpCatScores.Range("F1").Value = 11
$mlpCatScores.Range("G1").Value = 12
$mlpCatScores.Range("H1").Value = 13
$mlpCatScores.Range("I1").Value = 14
$mlpCatScores.Range("J1").Value = 15
$mlpCatScores.Range("K1").Value = "Mean"

$mlpCatScores.Range("A2").Value = "MAPE"
$mlpCatScores.Range("B2").Value = 3.09
$mlpCatScores.Range("C2").Value = 3.3
$mlpCatScores.Range("D2").Value = 3.43
$mlpCatScores.Range("E2").Value = 3.25
$mlpCatScores.Range("F2").Value = 3.09
$mlpCatScores.Range("G2").Value = 2.89
$mlpCatScores.Range("H2").Value = 3.08
$mlpCatScores.Range("I2").Value = 2.76
$mlpCatScores.Range("J2").Value = 2.76
$mlpCatScores.Range("K2").Value = 3.03

$mlpCatScores.Range("A3").Value = "RMSE"
$mlpCatScores.Range("B3").Value = 0.02
$mlpCatScores.Range("C3").Value = 0.02
$mlpCatScores.Range("D3").Value = 0.02
$mlpCatScores.Range("E3").Value = 0.01
$mlpCatScores.Range("F3").Value = 0.01
$mlpCatScores.Range("G3").Value = 0.01
$mlpCatScores.Range("H3").Value = 0.01
$mlpCatScores.Range("I3").Value = 0.01
$mlpCatScores.Range("J3").Value = 0.01
$mlpCatScores.Range("K3").Value = 0.01

$styleDonor.Copy() | Out-Null
$mlpCatScores.Range("B1:K1").PasteSpecial($xlPasteFormats) | Out-Null
$mlpCatScores.Range("A2").PasteSpecial($xlPasteFormats) | Out-Null
$mlpCatScores.Range("A3").PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------------
# 5) new ws_cat_scores sheet (sheetId 5), appended at the very end
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCatScores = $wb.Worksheets.Add($null, $lastSheet)
$wsCatScores.Name = "ws_cat_scores"

$wsCatScores.Range("B1").Value = 7
$wsCatScores.Range("C1").Value = 8
$wsCatScores.Range("D1").Value = 9
$wsCatScores.Range("E1").Value = 10
$wsCatScores.Range("F1").Value = 11
$wsCatScores.Range("G1").Value = 12
$wsCatScores.Range("H1").Value = 13
$wsCatScores.Range("I1").Value = 14
$wsCatScores.Range("J1").Value = 15
$wsCatScores.Range("K1").Value = "Mean"

$wsCatScores.Range("A2").Value = "MAPE"
$wsCatScores.Range("B2").Value = 23.49
$wsCatScores.Range("C2").Value = 22.62
$wsCatScores.Range("D2").Value = 22.39
$wsCatScores.Range("E2").Value = 21.17
$wsCatScores.Range("F2").Value = 21.68
$wsCatScores.Range("G2").Value = 22.01
$wsCatScores.Range("H2").Value = 21.97
$wsCatScores.Range("I2").Value = 21.83
$wsCatScores.Range("J2").Value = 21.72
$wsCatScores.Range("K2").Value = 22

$wsCatScores.Range("A3").Value = "RMSE"
$wsCatScores.Range("B3").Value = 0.1
$wsCatScores.Range("C3").Value = 0.1
$wsCatScores.Range("D3").Value = 0.09
$wsCatScores.Range("E3").Value = 0.09
$wsCatScores.Range("F3").Value = 0.09
$wsCatScores.Range("G3").Value = 0.09
$wsCatScores.Range("H3").Value = 0.09
$wsCatScores.Range("I3").Value = 0.08
$wsCatScores.Range("J3").Value = 0.08
$wsCatScores.Range("K3").Value = 0.09

$styleDonor.Copy() | Out-Null
$wsCatScores.Range("B1:K1").PasteSpecial($xlPasteFormats) | Out-Null
$wsCatScores.Range("A2").PasteSpecial($xlPasteFormats) | Out-Null
$wsCatScores.Range("A3").PasteSpecial($xlPasteFormats) | Out-Null

$runTimes.Select() | Out-Null
$runTimes.Range("A1").Select() | Out-Null
